# feat: add 2022-Q3 data
#
# 1. "总计" sheet: the existing row 2 (2021-Q4 totals) is pushed down to row 3,
#    and row 2 is overwritten with the new 2022-Q3 totals.
# 2. A brand-new sheet is appended at the end of the workbook, populated with
#    the (unchanged) fund-level detail that used to live on the "2021-Q4"
#    sheet, and named "2021-Q4".
# 3. The original "2021-Q4" sheet (now freed up) is renamed to "2022-Q3" and
#    its fund-level detail is replaced with the new quarter's numbers.
#
# NOTE: worksheet handles in this host resolve by *position*, so a variable
# captured before an Add/Move/Rename can silently point at the wrong sheet
# afterwards. Every lookup below is therefore done fresh, right before it is
# used, instead of being cached across a mutating call.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet - insert the new totals row, push the old one down.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("总计").Range("A2:D2").Copy($wb.Worksheets.Item("总计").Range("A3:D3"))
$wb.Worksheets.Item("总计").Cells.Item(3, 1).Value = 1

$wb.Worksheets.Item("总计").Cells.Item(2, 2).Value = "2022-Q3"
$wb.Worksheets.Item("总计").Cells.Item(2, 3).Value = 9
$wb.Worksheets.Item("总计").Cells.Item(2, 4).Value = 0.72

# ---------------------------------------------------------------------------
# 2) Preserve the old "2021-Q4" fund detail under a new tab at the end. A
#    full worksheet Copy (rather than Worksheets.Add + range copy) keeps the
#    sheet-level properties (sheetPr/sheetFormatPr/pageMargins) identical to
#    the original sheet instead of reverting to "blank sheet" defaults.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q4").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Free up the "2021-Q4" name from the original sheet (position 2) before
# claiming it on the freshly-appended copy ("2021-Q4 (2)", now the last
# sheet).
$wb.Worksheets.Item(2).Name = "2022-Q3"
$wb.Worksheets.Item("2021-Q4 (2)").Name = "2021-Q4"

# ---------------------------------------------------------------------------
# 3) Replace the "2022-Q3" sheet's fund data with the new quarter's numbers.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

# Extend the index column (A) style/formatting down to the 6 new rows.
$q3.Range("A4").Copy($q3.Range("A5:A10"))

# Force columns B-G to be stored as literal text (fund codes have leading
# zeros, and the scraped figures are text, not numbers) before writing them.
$q3.Range("B2:G10").NumberFormat = "@"

$q3Data = @(
    @(0, "011738", "华安兴安优选一年持有期混合A", "13.36", "55.44", "1.86", "0.2485", 6),
    @(1, "005695", "华安睿明两年定期开放灵活配置混合A", "4.27", "93.55", "4.84", "0.2067", 3),
    @(2, "011739", "华安兴安优选一年持有期混合C", "8.27", "55.44", "1.86", "0.1538", 6),
    @(3, "011390", "华安添祥6个月持有期混合A", "6.67", "33.77", "1.24", "0.0827", 6),
    @(4, "009409", "华安添福18个月持有期混合A", "0.72", "21.56", "2.18", "0.0157", 1),
    @(5, "003659", "山西证券策略精选灵活配置混合", "0.27", "78.35", "2.65", "0.0072", 10),
    @(6, "005696", "华安睿明两年定期开放灵活配置混合C", "0.07", "93.55", "4.84", "0.0034", 3),
    @(7, "009410", "华安添福18个月持有期混合C", "0.09", "21.56", "2.18", "0.0020", 1),
    @(8, "016181", "华安添祥6个月持有期混合C", "0.00", "33.77", "1.24", "0.00", 6)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# The very last market-value cell (G10) is stored as a genuine 0, not text.
$q3.Range("G10").NumberFormat = "General"
$q3.Cells.Item(10, 7).Value = 0
